# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.403.22"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "'3.772.59"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'615.64"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "'177.57"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "'3.767.91"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").Value = "'6.49"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "'39.85"
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").Value = "'4.401.43"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "'3.774.86"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "'69.463.51"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("E19").Value = "  -3.40%  "
$ws.Range("D20").Value = "'508.12"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'16.42"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").Value = "'9.39"
$ws.Range("E22").Value = "  -1.97%  "
$ws.Range("D23").Value = "'0.731"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").Value = "'2.48"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "'86.15"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("E26").Value = "  +5.65%  "
$ws.Range("D27").Value = "'12.87"
$ws.Range("D28").Value = "'10.51"
$ws.Range("E28").Value = "  -5.89%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'2.55"
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("D31").Value = "'3.02"
$ws.Range("E31").Value = "  +3.59%  "
$ws.Range("D32").Value = "'8.11"
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("D33").Value = "'31.03"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").Value = "'6.12"
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("E38").Value = "  +6.25%  "
$ws.Range("D39").Value = "'0.340"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").Value = "'464.77"
$ws.Range("E40").Value = "  +9.72%  "
$ws.Range("D42").Value = "'3.02"
$ws.Range("E42").Value = "  +9.03%  "
$ws.Range("D43").Value = "'49.83"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("D44").Value = "'44.30"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").Value = "'8.59"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").Value = "'2.950.49"
$ws.Range("D47").Value = "'0.0362"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").Value = "'27.29"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D50").Value = "'139.16"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "'2.46"
$ws.Range("E51").Value = "  -0.78%  "
